$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.796.57'
$ws.Range("E2").Value = '  +2.90%  '

$ws.Range("D3").Value = '1.881.34'
$ws.Range("E3").Value = '  +3.10%  '

$ws.Range("E4").Value = '  +0.67%  '

$ws.Range("D5").Value = "'324.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").Value = "'0.4678"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.12%  '

$ws.Range("D8").Value = "'0.3934"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.36%  '

$ws.Range("D9").Value = "'0.07933"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("D10").Value = "'0.9787"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.05%  '

$ws.Range("E11").Value = '  +2.12%  '

$ws.Range("D12").Value = '1.847.49'
$ws.Range("E12").Value = '  +2.86%  '

$ws.Range("D13").Value = "'7.020"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.38%  '

$ws.Range("D14").Value = "'5.740"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.67%  '

$ws.Range("D15").Value = "'0.06962"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.48%  '

$ws.Range("D16").Value = "'88.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.71%  '

$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").Value = "'0.00001010"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.58%  '

$ws.Range("D19").Value = "'16.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.90%  '

$ws.Range("E20").Value = '  +0.63%  '

$ws.Range("D21").Value = '28.817.85'
$ws.Range("E21").Value = '  +2.90%  '

$ws.Range("D22").Value = "'5.346"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("D23").Value = "'11.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.23%  '

$ws.Range("E24").Value = '  +1.31%  '

$ws.Range("D25").Value = '2.139.59'
$ws.Range("E25").Value = '  +6.01%  '

$ws.Range("D26").Value = "'153.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.99%  '

$ws.Range("D27").Value = "'19.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '

$ws.Range("D28").Value = "'5.756"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.54%  '

$ws.Range("D29").Value = "'2.000"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.37%  '

$ws.Range("D30").Value = "'120.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.95%  '

$ws.Range("D31").Value = "'0.09397"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.84%  '

$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("D33").Value = "'5.317"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.62%  '

$ws.Range("D34").Value = "'1.356"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.24%  '

$ws.Range("D35").Value = "'3.353"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.38%  '

$ws.Range("D36").Value = "'0.05917"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").Value = "'0.02124"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.66%  '

$ws.Range("E38").Value = '  +1.43%  '

$ws.Range("E39").Value = '  +4.53%  '

$ws.Range("D40").Value = "'0.5721"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.68%  '

$ws.Range("D41").Value = "'0.1797"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.96%  '

$ws.Range("D42").Value = "'9.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("D43").Value = "'0.07331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.78%  '

$ws.Range("D44").Value = "'11.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.33%  '

$ws.Range("D45").Value = "'0.5347"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.01%  '

$ws.Range("D46").Value = "'1.153"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.39%  '

$ws.Range("D47").Value = "'1.847"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.45%  '

$ws.Range("D48").Value = "'2.106"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.49%  '

$ws.Range("E49").Value = '  +1.92%  '

$ws.Range("D50").Value = "'2.372"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.82%  '

$ws.Range("D51").Value = "'1.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.58%  '
